$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 402542.88
$ws.Range("J17").Value = 402542.88
$ws.Range("L17").Value = 1207628.64
$ws.Range("N17").Value = -1207964.64

# Row 33
$ws.Range("H33").Value = 860.4231
$ws.Range("I33").Value = 562.7646999999999
$ws.Range("K33").Value = 562.7646999999999
$ws.Range("M33").Value = -333.7646999999999

# Row 76
$ws.Range("H76").Value = 2835.4443
$ws.Range("I76").Value = 2665.5334
$ws.Range("J76").Value = 3685
$ws.Range("K76").Value = 2665.5334
$ws.Range("L76").Value = 3685
$ws.Range("M76").Value = -2350.5334
$ws.Range("N76").Value = -4315

# Row 79
$ws.Range("H79").Value = 2835.4443
$ws.Range("I79").Value = 2665.5334
$ws.Range("J79").Value = 3685
$ws.Range("K79").Value = 2665.5334
$ws.Range("L79").Value = 3685
$ws.Range("M79").Value = -1573.5334
$ws.Range("N79").Value = -5869

# Row 87
$ws.Range("H87").Value = 29376.695
$ws.Range("J87").Value = 29376.695
$ws.Range("L87").Value = 29376.695
$ws.Range("N87").Value = -31872.695

# Row 90
$ws.Range("H90").Value = 29376.695
$ws.Range("J90").Value = 29376.695
$ws.Range("L90").Value = 88130.08499999999
$ws.Range("N90").Value = -100610.085

# Row 92
$ws.Range("H92").Value = 336.58334
$ws.Range("I92").Value = 368.51852
$ws.Range("J92").Value = 240.77777
$ws.Range("K92").Value = 368.51852
$ws.Range("L92").Value = 240.77777
$ws.Range("M92").Value = 879.4814799999999
$ws.Range("N92").Value = -2736.77777

# Row 106
$ws.Range("H106").Value = 3407.6924
$ws.Range("I106").Value = 2216.6667
$ws.Range("J106").Value = 4428.5713
$ws.Range("K106").Value = 2216.6667
$ws.Range("L106").Value = 4428.5713
$ws.Range("M106").Value = -1585.6667
$ws.Range("N106").Value = -5690.5713

# Row 111
$ws.Range("H111").Value = 1300
$ws.Range("I111").Value = 1350
$ws.Range("K111").Value = 4050
$ws.Range("M111").Value = -983

# Row 113
$ws.Range("H113").Value = 6040.6
$ws.Range("I113").Value = 4485.7144
$ws.Range("J113").Value = 9668.666999999999
$ws.Range("K113").Value = 4485.7144
$ws.Range("L113").Value = 9668.666999999999
$ws.Range("M113").Value = -1231.7144
$ws.Range("N113").Value = -16176.667

# Row 132
$ws.Range("H132").Value = 2633119.5
$ws.Range("I132").Value = 2741097.8
$ws.Range("K132").Value = 8223293.399999999
$ws.Range("M132").Value = -8220763.399999999

# Row 137
$ws.Range("H137").Value = 1821446.9
$ws.Range("I137").Value = 2706310.8
$ws.Range("J137").Value = 2559.7222
$ws.Range("K137").Value = 8118932.399999999
$ws.Range("L137").Value = 7679.1666
$ws.Range("M137").Value = -8116382.399999999
$ws.Range("N137").Value = -12779.1666

$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 18999.857
$ws.Range("I33").Value = 18999
$ws.Range("K33").Value = 18999
$ws.Range("M33").Value = -18670

# Row 52
$ws.Range("H52").Value = 27500
$ws.Range("J52").Value = 27500
$ws.Range("L52").Value = 27500
$ws.Range("N52").Value = -28136

# Row 132
$ws.Range("H132").Value = 2516.3462
$ws.Range("I132").Value = 1886.5294
$ws.Range("J132").Value = 3706
$ws.Range("K132").Value = 5659.5882
$ws.Range("L132").Value = 11118
$ws.Range("M132").Value = -3129.5882
$ws.Range("N132").Value = -16178

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 23626.715
$ws.Range("J82").Value = 28743.182
$ws.Range("L82").Value = 28743.182
$ws.Range("N82").Value = -29509.182

# Row 85
$ws.Range("H85").Value = 23626.715
$ws.Range("J85").Value = 28743.182
$ws.Range("L85").Value = 28743.182
$ws.Range("N85").Value = -31395.182

# Row 94
$ws.Range("H94").Value = 949.4
$ws.Range("I94").Value = 789.8570999999999
$ws.Range("J94").Value = 1321.6666
$ws.Range("K94").Value = 789.8570999999999
$ws.Range("L94").Value = 1321.6666
$ws.Range("M94").Value = -338.8570999999999
$ws.Range("N94").Value = -2223.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1278.0857
$ws.Range("I107").Value = 1152.8
$ws.Range("K107").Value = 1152.8
$ws.Range("M107").Value = 767.2

# Row 132
$ws.Range("H132").Value = 2393.5945
$ws.Range("I132").Value = 1784.6364
$ws.Range("J132").Value = 3286.7334
$ws.Range("K132").Value = 5353.9092
$ws.Range("L132").Value = 9860.200199999999
$ws.Range("M132").Value = -2823.9092
$ws.Range("N132").Value = -14920.2002

# Row 135
$ws.Range("H135").Value = 26698.75
$ws.Range("J135").Value = 26698.75
$ws.Range("L135").Value = 26698.75
$ws.Range("N135").Value = -36838.75

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 1000.8889
$ws.Range("I107").Value = 454
$ws.Range("J107").Value = 1438.4
$ws.Range("K107").Value = 454
$ws.Range("L107").Value = 1438.4
$ws.Range("M107").Value = 1466
$ws.Range("N107").Value = -5278.4

# Row 126
$ws.Range("H126").Value = 629387
$ws.Range("I126").Value = 3284
$ws.Range("K126").Value = 9852
$ws.Range("M126").Value = -7382

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 359643
$ws.Range("J2").Value = 8750.5
$ws.Range("L2").Value = 8750.5
$ws.Range("N2").Value = -8974.5

# Row 14
$ws.Range("H14").Value = 178816.5
$ws.Range("I14").Value = 1255002
$ws.Range("J14").Value = 25075.715
$ws.Range("K14").Value = 1255002
$ws.Range("L14").Value = 25075.715
$ws.Range("M14").Value = -1254830
$ws.Range("N14").Value = -25419.715

# Row 46
$ws.Range("H46").Value = 1479.2307
$ws.Range("I46").Value = 769.1667
$ws.Range("K46").Value = 769.1667
$ws.Range("M46").Value = -581.1667

# Row 61
$ws.Range("H61").Value = 125004580
$ws.Range("I61").Value = 200000770
$ws.Range("J61").Value = 10933.333
$ws.Range("K61").Value = 200000770
$ws.Range("L61").Value = 10933.333
$ws.Range("M61").Value = -200000568
$ws.Range("N61").Value = -11337.333

# Row 93
$ws.Range("H93").Value = 2208.0908
$ws.Range("I93").Value = 1634.4286
$ws.Range("J93").Value = 3212
$ws.Range("K93").Value = 1634.4286
$ws.Range("L93").Value = 3212
$ws.Range("M93").Value = -386.4286
$ws.Range("N93").Value = -5708

# Row 106
$ws.Range("H106").Value = 25795.875
$ws.Range("J106").Value = 25795.875
$ws.Range("L106").Value = 25795.875
$ws.Range("N106").Value = -28319.875

# Row 109
$ws.Range("H109").Value = 14753
$ws.Range("I109").Value = 10259
$ws.Range("J109").Value = 17000
$ws.Range("K109").Value = 10259
$ws.Range("L109").Value = 17000
$ws.Range("M109").Value = -8872
$ws.Range("N109").Value = -19774

# Row 113
$ws.Range("H113").Value = 125004580
$ws.Range("I113").Value = 200000770
$ws.Range("J113").Value = 10933.333
$ws.Range("K113").Value = 200000770
$ws.Range("L113").Value = 10933.333
$ws.Range("M113").Value = -199998600
$ws.Range("N113").Value = -15273.333

# Row 122
$ws.Range("H122").Value = 3891.25
$ws.Range("I122").Value = 2862.2222
$ws.Range("J122").Value = 5214.2856
$ws.Range("K122").Value = 8586.6666
$ws.Range("L122").Value = 15642.8568
$ws.Range("M122").Value = -6136.6666
$ws.Range("N122").Value = -20542.8568

# Row 125
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840

# Row 132
$ws.Range("H132").Value = 2773
$ws.Range("I132").Value = 2028.375
$ws.Range("J132").Value = 3765.8333
$ws.Range("K132").Value = 6085.125
$ws.Range("L132").Value = 11297.4999
$ws.Range("M132").Value = -3555.125
$ws.Range("N132").Value = -16357.4999

# Row 136
$ws.Range("H136").Value = 2131769
$ws.Range("I136").Value = 3336351
$ws.Range("J136").Value = 6035.8823
$ws.Range("K136").Value = 10009053
$ws.Range("L136").Value = 18107.6469
$ws.Range("M136").Value = -10006503
$ws.Range("N136").Value = -23207.6469

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 113
$ws.Range("H113").Value = 2675
$ws.Range("I113").Value = 200
$ws.Range("J113").Value = 5150
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 15450
$ws.Range("M113").Value = 1570
$ws.Range("N113").Value = -19790

# Row 132
$ws.Range("H132").Value = 1571663.2
$ws.Range("I132").Value = 1925833.5
$ws.Range("J132").Value = 36925.668
$ws.Range("K132").Value = 5777500.5
$ws.Range("L132").Value = 110777.004
$ws.Range("M132").Value = -5774970.5
$ws.Range("N132").Value = -115837.004

# Row 139
$ws.Range("H139").Value = 54900
$ws.Range("J139").Value = 54900
$ws.Range("L139").Value = 54900
$ws.Range("N139").Value = -65180

# Row 141
$ws.Range("H141").Value = 28681.818
$ws.Range("J141").Value = 28681.818
$ws.Range("L141").Value = 28681.818
$ws.Range("N141").Value = -39041.818
